$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1631.7142
$ws.Range("I2").Value = 1081.8
$ws.Range("K2").Value = 1081.8
$ws.Range("M2").Value = -968.8
$ws.Range("H32").Value = 5398.946
$ws.Range("I32").Value = 3698.3767
$ws.Range("J32").Value = 14128.533
$ws.Range("K32").Value = 3698.3767
$ws.Range("L32").Value = 14128.533
$ws.Range("M32").Value = -3411.3767
$ws.Range("N32").Value = -14702.533
$ws.Range("H45").Value = 6077.1055
$ws.Range("I45").Value = 7035.3125
$ws.Range("J45").Value = 966.6667
$ws.Range("K45").Value = 7035.3125
$ws.Range("L45").Value = 966.6667
$ws.Range("M45").Value = -6658.3125
$ws.Range("N45").Value = -1720.6667
$ws.Range("H61").Value = 3920.5
$ws.Range("I61").Value = 4572.533
$ws.Range("J61").Value = 1964.4
$ws.Range("K61").Value = 4572.533
$ws.Range("L61").Value = 1964.4
$ws.Range("M61").Value = -4360.533
$ws.Range("N61").Value = -2388.4
$ws.Range("H116").Value = 1631.7142
$ws.Range("I116").Value = 1081.8
$ws.Range("K116").Value = 1081.8
$ws.Range("M116").Value = 1212.2
$ws.Range("H122").Value = 1834479.4
$ws.Range("I122").Value = 3667976.8
$ws.Range("J122").Value = 982.1429000000001
$ws.Range("K122").Value = 11003930.4
$ws.Range("L122").Value = 2946.4287
$ws.Range("M122").Value = -11001480.4
$ws.Range("N122").Value = -7846.4287
$ws.Range("H132").Value = 1540505.6
$ws.Range("I132").Value = 1450.2941
$ws.Range("J132").Value = 7147064.5
$ws.Range("K132").Value = 4350.8823
$ws.Range("L132").Value = 21441193.5
$ws.Range("M132").Value = -1820.8823
$ws.Range("N132").Value = -21446253.5
$ws.Range("H136").Value = 3920.5
$ws.Range("I136").Value = 4572.533
$ws.Range("J136").Value = 1964.4
$ws.Range("K136").Value = 13717.599
$ws.Range("L136").Value = 5893.200000000001
$ws.Range("M136").Value = -11167.599
$ws.Range("N136").Value = -10993.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1631.7142
$ws.Range("I3").Value = 1081.8
$ws.Range("K3").Value = 1081.8
$ws.Range("M3").Value = -967.8
$ws.Range("H82").Value = 17006.9
$ws.Range("I82").Value = 3317.5
$ws.Range("J82").Value = 26133.166
$ws.Range("K82").Value = 3317.5
$ws.Range("L82").Value = 26133.166
$ws.Range("M82").Value = -2934.5
$ws.Range("N82").Value = -26899.166
$ws.Range("H85").Value = 17006.9
$ws.Range("I85").Value = 3317.5
$ws.Range("J85").Value = 26133.166
$ws.Range("K85").Value = 3317.5
$ws.Range("L85").Value = 26133.166
$ws.Range("M85").Value = -1991.5
$ws.Range("N85").Value = -28785.166
$ws.Range("H86").Value = 1847.091
$ws.Range("I86").Value = 1789.75
$ws.Range("K86").Value = 1789.75
$ws.Range("M86").Value = -666.75
$ws.Range("H89").Value = 1847.091
$ws.Range("I89").Value = 1789.75
$ws.Range("K89").Value = 8948.75
$ws.Range("M89").Value = -3332.75
$ws.Range("H134").Value = 3122.5356
$ws.Range("I134").Value = 3188.8125
$ws.Range("J134").Value = 2724.875
$ws.Range("K134").Value = 9566.4375
$ws.Range("L134").Value = 8174.625
$ws.Range("M134").Value = -7031.4375
$ws.Range("N134").Value = -13244.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5146.5
$ws.Range("I31").Value = 1293.1578
$ws.Range("J31").Value = 13281.333
$ws.Range("K31").Value = 1293.1578
$ws.Range("L31").Value = 13281.333
$ws.Range("M31").Value = -998.1578
$ws.Range("N31").Value = -13871.333
$ws.Range("H34").Value = 5146.5
$ws.Range("I34").Value = 1293.1578
$ws.Range("J34").Value = 13281.333
$ws.Range("K34").Value = 1293.1578
$ws.Range("L34").Value = 13281.333
$ws.Range("M34").Value = -1091.1578
$ws.Range("N34").Value = -13685.333
$ws.Range("H58").Value = 1260.5532
$ws.Range("I58").Value = 752.24243
$ws.Range("J58").Value = 2458.7144
$ws.Range("K58").Value = 752.24243
$ws.Range("L58").Value = 2458.7144
$ws.Range("M58").Value = -549.24243
$ws.Range("N58").Value = -2864.7144
$ws.Range("H62").Value = 6797
$ws.Range("I62").Value = 7496.25
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 7496.25
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -6872.25
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 6797
$ws.Range("I65").Value = 7496.25
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 37481.25
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -34361.25
$ws.Range("N65").Value = -26240
$ws.Range("H136").Value = 1260.5532
$ws.Range("I136").Value = 752.24243
$ws.Range("J136").Value = 2458.7144
$ws.Range("K136").Value = 2256.72729
$ws.Range("L136").Value = 7376.1432
$ws.Range("M136").Value = 293.2727100000002
$ws.Range("N136").Value = -12476.1432
$ws.Range("H141").Value = 34071.25
$ws.Range("J141").Value = 34071.25
$ws.Range("L141").Value = 34071.25
$ws.Range("N141").Value = -44431.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 102050.26
$ws.Range("J5").Value = 273366.97
$ws.Range("L5").Value = 820100.9099999999
$ws.Range("N5").Value = -820324.9099999999
$ws.Range("H92").Value = 417.5
$ws.Range("J92").Value = 421
$ws.Range("L92").Value = 1263
$ws.Range("N92").Value = -3759
$ws.Range("H117").Value = 22233352
$ws.Range("I117").Value = 17038.166
$ws.Range("J117").Value = 37044228
$ws.Range("K117").Value = 51114.49800000001
$ws.Range("L117").Value = 111132684
$ws.Range("M117").Value = -47672.49800000001
$ws.Range("N117").Value = -111139568
$ws.Range("H121").Value = 970.62067
$ws.Range("I121").Value = 496
$ws.Range("J121").Value = 1069.5
$ws.Range("K121").Value = 1488
$ws.Range("L121").Value = 3208.5
$ws.Range("M121").Value = -178
$ws.Range("N121").Value = -5828.5
$ws.Range("H129").Value = 1203
$ws.Range("I129").Value = 757.5
$ws.Range("J129").Value = 1648.5
$ws.Range("K129").Value = 2272.5
$ws.Range("L129").Value = 4945.5
$ws.Range("M129").Value = 2727.5
$ws.Range("N129").Value = -14945.5
$ws.Range("H135").Value = 102050.26
$ws.Range("J135").Value = 273366.97
$ws.Range("L135").Value = 2460302.73
$ws.Range("N135").Value = -2465372.73

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2402258.2
$ws.Range("I122").Value = 4052547.5
$ws.Range("J122").Value = 1837.5454
$ws.Range("K122").Value = 12157642.5
$ws.Range("L122").Value = 5512.6362
$ws.Range("M122").Value = -12155192.5
$ws.Range("N122").Value = -10412.6362
$ws.Range("H132").Value = 1668.6666
$ws.Range("I132").Value = 1193.5476
$ws.Range("J132").Value = 3331.5833
$ws.Range("K132").Value = 3580.642800000001
$ws.Range("L132").Value = 9994.749899999999
$ws.Range("M132").Value = -1050.642800000001
$ws.Range("N132").Value = -15054.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1882.1666
$ws.Range("I61").Value = 2122
$ws.Range("J61").Value = 1402.5
$ws.Range("K61").Value = 2122
$ws.Range("L61").Value = 1402.5
$ws.Range("M61").Value = -1920
$ws.Range("N61").Value = -1806.5
$ws.Range("H113").Value = 1882.1666
$ws.Range("I113").Value = 2122
$ws.Range("J113").Value = 1402.5
$ws.Range("K113").Value = 2122
$ws.Range("L113").Value = 1402.5
$ws.Range("M113").Value = 48
$ws.Range("N113").Value = -5742.5
$ws.Range("H122").Value = 2308547.5
$ws.Range("I122").Value = 2980994.5
$ws.Range("K122").Value = 8942983.5
$ws.Range("M122").Value = -8940533.5
$ws.Range("H136").Value = 4705.1665
$ws.Range("I136").Value = 4024.8838
$ws.Range("J136").Value = 7364.4546
$ws.Range("K136").Value = 12074.6514
$ws.Range("L136").Value = 22093.3638
$ws.Range("M136").Value = -9524.651400000001
$ws.Range("N136").Value = -27193.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1554.4445
$ws.Range("I122").Value = 1466.1666
$ws.Range("J122").Value = 1731
$ws.Range("K122").Value = 4398.4998
$ws.Range("L122").Value = 5193
$ws.Range("M122").Value = -1948.4998
$ws.Range("N122").Value = -10093
$ws.Range("H123").Value = 28971.428
$ws.Range("J123").Value = 28971.428
$ws.Range("L123").Value = 28971.428
$ws.Range("N123").Value = -38771.428
$ws.Range("H132").Value = 30499.912
$ws.Range("I132").Value = 46410.59
$ws.Range("J132").Value = 1330.3334
$ws.Range("K132").Value = 139231.77
$ws.Range("L132").Value = 3991.0002
$ws.Range("M132").Value = -136701.77
$ws.Range("N132").Value = -9051.0002
$ws.Range("H136").Value = 7814839
$ws.Range("I136").Value = 2626.195
$ws.Range("J136").Value = 21740956
$ws.Range("K136").Value = 7878.585000000001
$ws.Range("L136").Value = 65222868
$ws.Range("M136").Value = -5328.585000000001
$ws.Range("N136").Value = -65227968
